# Atualizacao de bases das ligas - swap da liga Estonia Meistriliiga
# Applies corrected match/team/odds data by swapping paired rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 32 and 33 (columns B:AC) ---
$ws.Range("B32").Value = 6475429
$ws.Range("C32").Value = "Estonia Meistriliiga"
$ws.Range("D32").Value = "Estonia Meistriliiga"
$ws.Range("E32").Value = 45035.5
$ws.Range("F32").Value = "Parnu JK Vaprus"
$ws.Range("G32").Value = "JK Trans Narva"
$ws.Range("H32").Value = 1
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = "H"
$ws.Range("K32").Value = 2.4
$ws.Range("L32").Value = 3.6
$ws.Range("M32").Value = 2.4
$ws.Range("N32").Value = 2.8
$ws.Range("O32").Value = 3.6
$ws.Range("P32").Value = 2.1
$ws.Range("Q32").Value = 0.25
$ws.Range("R32").Value = 1.9
$ws.Range("S32").Value = 1.9
$ws.Range("T32").Value = 2.5
$ws.Range("U32").Value = 1.925
$ws.Range("V32").Value = 1.875
$ws.Range("W32").Value = 1.8
$ws.Range("X32").Value = -1
$ws.Range("Y32").Value = -1
$ws.Range("Z32").Value = 0.8999999999999999
$ws.Range("AA32").Value = -1
$ws.Range("AB32").Value = -1
$ws.Range("AC32").Value = 0.875
$ws.Range("B33").Value = 6478314
$ws.Range("C33").Value = "Estonia Meistriliiga"
$ws.Range("D33").Value = "Estonia Meistriliiga"
$ws.Range("E33").Value = 45035.5
$ws.Range("F33").Value = "JK Tammeka Tartu"
$ws.Range("G33").Value = "JK Tallinna Kalev"
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = 2
$ws.Range("J33").Value = "A"
$ws.Range("K33").Value = 2.4
$ws.Range("L33").Value = 3.6
$ws.Range("M33").Value = 2.4
$ws.Range("N33").Value = 2.3
$ws.Range("O33").Value = 3.4
$ws.Range("P33").Value = 2.6
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = 1.725
$ws.Range("S33").Value = 1.975
$ws.Range("T33").Value = 2.25
$ws.Range("U33").Value = 1.9
$ws.Range("V33").Value = 1.9
$ws.Range("W33").Value = -1
$ws.Range("X33").Value = -1
$ws.Range("Y33").Value = 1.6
$ws.Range("Z33").Value = -1
$ws.Range("AA33").Value = 0.9750000000000001
$ws.Range("AB33").Value = 0.8999999999999999
$ws.Range("AC33").Value = -1

# --- Swap rows 77 and 78 (columns B:AC) ---
$ws.Range("B77").Value = 6139017
$ws.Range("C77").Value = "Estonia Meistriliiga"
$ws.Range("D77").Value = "Estonia Meistriliiga"
$ws.Range("E77").Value = 45084.5
$ws.Range("F77").Value = "JK Tammeka Tartu"
$ws.Range("G77").Value = "Harju JK Laagri"
$ws.Range("H77").Value = 2
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = "H"
$ws.Range("K77").Value = 1.666
$ws.Range("L77").Value = 3.6
$ws.Range("M77").Value = 4.2
$ws.Range("N77").Value = 1.727
$ws.Range("O77").Value = 3.5
$ws.Range("P77").Value = 4
$ws.Range("Q77").Value = -0.75
$ws.Range("R77").Value = 2
$ws.Range("S77").Value = 1.8
$ws.Range("T77").Value = 2.5
$ws.Range("U77").Value = 1.9
$ws.Range("V77").Value = 1.9
$ws.Range("W77").Value = 0.7270000000000001
$ws.Range("X77").Value = -1
$ws.Range("Y77").Value = -1
$ws.Range("Z77").Value = 1
$ws.Range("AA77").Value = -1
$ws.Range("AB77").Value = -1
$ws.Range("AC77").Value = 0.8999999999999999
$ws.Range("B78").Value = 6139018
$ws.Range("C78").Value = "Estonia Meistriliiga"
$ws.Range("D78").Value = "Estonia Meistriliiga"
$ws.Range("E78").Value = 45084.5
$ws.Range("F78").Value = "JK Tallinna Kalev"
$ws.Range("G78").Value = "JK Trans Narva"
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 1
$ws.Range("J78").Value = "A"
$ws.Range("K78").Value = 2.4
$ws.Range("L78").Value = 3.4
$ws.Range("M78").Value = 2.5
$ws.Range("N78").Value = 2.875
$ws.Range("O78").Value = 3.1
$ws.Range("P78").Value = 2.3
$ws.Range("Q78").Value = 0.25
$ws.Range("R78").Value = 1.75
$ws.Range("S78").Value = 2.05
$ws.Range("T78").Value = 2.25
$ws.Range("U78").Value = 1.925
$ws.Range("V78").Value = 1.875
$ws.Range("W78").Value = -1
$ws.Range("X78").Value = -1
$ws.Range("Y78").Value = 1.3
$ws.Range("Z78").Value = -1
$ws.Range("AA78").Value = 1.05
$ws.Range("AB78").Value = -1
$ws.Range("AC78").Value = 0.875

# --- Swap rows 137 and 138 (columns B:AC) ---
$ws.Range("B137").Value = 6139064
$ws.Range("C137").Value = "Estonia Meistriliiga"
$ws.Range("D137").Value = "Estonia Meistriliiga"
$ws.Range("E137").Value = 45192.35416666666
$ws.Range("F137").Value = "JK Trans Narva"
$ws.Range("G137").Value = "Harju JK Laagri"
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 3
$ws.Range("J137").Value = "A"
$ws.Range("K137").Value = 1.75
$ws.Range("L137").Value = 3.6
$ws.Range("M137").Value = 3.8
$ws.Range("N137").Value = 1.45
$ws.Range("O137").Value = 4
$ws.Range("P137").Value = 6
$ws.Range("Q137").Value = -1
$ws.Range("R137").Value = 1.85
$ws.Range("S137").Value = 1.95
$ws.Range("T137").Value = 2.5
$ws.Range("U137").Value = 1.9
$ws.Range("V137").Value = 1.9
$ws.Range("W137").Value = -1
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = 5
$ws.Range("Z137").Value = -1
$ws.Range("AA137").Value = 0.95
$ws.Range("AB137").Value = 0.8999999999999999
$ws.Range("AC137").Value = -1
$ws.Range("B138").Value = 6139067
$ws.Range("C138").Value = "Estonia Meistriliiga"
$ws.Range("D138").Value = "Estonia Meistriliiga"
$ws.Range("E138").Value = 45192.35416666666
$ws.Range("F138").Value = "Paide Linnameeskond"
$ws.Range("G138").Value = "Parnu JK Vaprus"
$ws.Range("H138").Value = 3
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = "H"
$ws.Range("K138").Value = 1.8
$ws.Range("L138").Value = 3.4
$ws.Range("M138").Value = 3.8
$ws.Range("N138").Value = 1.5
$ws.Range("O138").Value = 3.8
$ws.Range("P138").Value = 6
$ws.Range("Q138").Value = -1
$ws.Range("R138").Value = 1.75
$ws.Range("S138").Value = 1.95
$ws.Range("T138").Value = 2.5
$ws.Range("U138").Value = 1.9
$ws.Range("V138").Value = 1.9
$ws.Range("W138").Value = 0.5
$ws.Range("X138").Value = -1
$ws.Range("Y138").Value = -1
$ws.Range("Z138").Value = 0.75
$ws.Range("AA138").Value = -1
$ws.Range("AB138").Value = 0.8999999999999999
$ws.Range("AC138").Value = -1

# --- Swap rows 168 and 169 (columns B:AC) ---
$ws.Range("B168").Value = 6482819
$ws.Range("C168").Value = "Estonia Meistriliiga"
$ws.Range("D168").Value = "Estonia Meistriliiga"
$ws.Range("E168").Value = 45231.54166666666
$ws.Range("F168").Value = "JK Tammeka Tartu"
$ws.Range("G168").Value = "FC Kuressaare"
$ws.Range("H168").Value = 0
$ws.Range("I168").Value = 1
$ws.Range("J168").Value = "A"
$ws.Range("K168").Value = 1.833
$ws.Range("L168").Value = 3.5
$ws.Range("M168").Value = 3.5
$ws.Range("N168").Value = 2.1
$ws.Range("O168").Value = 3.4
$ws.Range("P168").Value = 2.875
$ws.Range("Q168").Value = -0.25
$ws.Range("R168").Value = 1.975
$ws.Range("S168").Value = 1.825
$ws.Range("T168").Value = 3
$ws.Range("U168").Value = 1.825
$ws.Range("V168").Value = 1.975
$ws.Range("W168").Value = -1
$ws.Range("X168").Value = -1
$ws.Range("Y168").Value = 1.875
$ws.Range("Z168").Value = -1
$ws.Range("AA168").Value = 0.825
$ws.Range("AB168").Value = -1
$ws.Range("AC168").Value = 0.9750000000000001
$ws.Range("B169").Value = 6416370
$ws.Range("C169").Value = "Estonia Meistriliiga"
$ws.Range("D169").Value = "Estonia Meistriliiga"
$ws.Range("E169").Value = 45231.54166666666
$ws.Range("F169").Value = "FC Levadia Tallinn"
$ws.Range("G169").Value = "Parnu JK Vaprus"
$ws.Range("H169").Value = 0
$ws.Range("I169").Value = 0
$ws.Range("J169").Value = "D"
$ws.Range("K169").Value = 1.166
$ws.Range("L169").Value = 7
$ws.Range("M169").Value = 11
$ws.Range("N169").Value = 1.2
$ws.Range("O169").Value = 6
$ws.Range("P169").Value = 11
$ws.Range("Q169").Value = -2
$ws.Range("R169").Value = 1.85
$ws.Range("S169").Value = 1.95
$ws.Range("T169").Value = 3
$ws.Range("U169").Value = 1.85
$ws.Range("V169").Value = 1.95
$ws.Range("W169").Value = -1
$ws.Range("X169").Value = 5
$ws.Range("Y169").Value = -1
$ws.Range("Z169").Value = -1
$ws.Range("AA169").Value = 0.95
$ws.Range("AB169").Value = -1
$ws.Range("AC169").Value = 0.95

# --- Swap rows 177 and 179 (columns B:AC) ---
$ws.Range("B177").Value = 6537957
$ws.Range("C177").Value = "Estonia Meistriliiga"
$ws.Range("D177").Value = "Estonia Meistriliiga"
$ws.Range("E177").Value = 45241.375
$ws.Range("F177").Value = "FC Flora Tallinn"
$ws.Range("G177").Value = "JK Nomme Kalju"
$ws.Range("H177").Value = 0
$ws.Range("I177").Value = 0
$ws.Range("J177").Value = "D"
$ws.Range("K177").Value = 1.4
$ws.Range("L177").Value = 4
$ws.Range("M177").Value = 7.5
$ws.Range("N177").Value = 1.5
$ws.Range("O177").Value = 4.2
$ws.Range("P177").Value = 5
$ws.Range("Q177").Value = -1
$ws.Range("R177").Value = 1.85
$ws.Range("S177").Value = 1.95
$ws.Range("T177").Value = 2.75
$ws.Range("U177").Value = 1.85
$ws.Range("V177").Value = 1.95
$ws.Range("W177").Value = -1
$ws.Range("X177").Value = 3.2
$ws.Range("Y177").Value = -1
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = 0.95
$ws.Range("AB177").Value = -1
$ws.Range("AC177").Value = 0.95
$ws.Range("B179").Value = 6533597
$ws.Range("C179").Value = "Estonia Meistriliiga"
$ws.Range("D179").Value = "Estonia Meistriliiga"
$ws.Range("E179").Value = 45241.375
$ws.Range("F179").Value = "FC Kuressaare"
$ws.Range("G179").Value = "Parnu JK Vaprus"
$ws.Range("H179").Value = 1
$ws.Range("I179").Value = 0
$ws.Range("J179").Value = "H"
$ws.Range("K179").Value = 2.5
$ws.Range("L179").Value = 3.4
$ws.Range("M179").Value = 2.5
$ws.Range("N179").Value = 2.15
$ws.Range("O179").Value = 3.6
$ws.Range("P179").Value = 2.875
$ws.Range("Q179").Value = -0.25
$ws.Range("R179").Value = 1.95
$ws.Range("S179").Value = 1.85
$ws.Range("T179").Value = 2.75
$ws.Range("U179").Value = 1.95
$ws.Range("V179").Value = 1.85
$ws.Range("W179").Value = 1.15
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = -1
$ws.Range("Z179").Value = 0.95
$ws.Range("AA179").Value = -1
$ws.Range("AB179").Value = -1
$ws.Range("AC179").Value = 0.8500000000000001

# --- Swap rows 178 and 180 (columns B:AC) ---
$ws.Range("B178").Value = 6535416
$ws.Range("C178").Value = "Estonia Meistriliiga"
$ws.Range("D178").Value = "Estonia Meistriliiga"
$ws.Range("E178").Value = 45241.375
$ws.Range("F178").Value = "Paide Linnameeskond"
$ws.Range("G178").Value = "FC Levadia Tallinn"
$ws.Range("H178").Value = 2
$ws.Range("I178").Value = 2
$ws.Range("J178").Value = "D"
$ws.Range("K178").Value = 3
$ws.Range("L178").Value = 3.8
$ws.Range("M178").Value = 2
$ws.Range("N178").Value = 3
$ws.Range("O178").Value = 4
$ws.Range("P178").Value = 1.909
$ws.Range("Q178").Value = 0.5
$ws.Range("R178").Value = 1.85
$ws.Range("S178").Value = 1.95
$ws.Range("T178").Value = 2.75
$ws.Range("U178").Value = 1.95
$ws.Range("V178").Value = 1.85
$ws.Range("W178").Value = -1
$ws.Range("X178").Value = 3
$ws.Range("Y178").Value = -1
$ws.Range("Z178").Value = 0.8500000000000001
$ws.Range("AA178").Value = -1
$ws.Range("AB178").Value = 0.95
$ws.Range("AC178").Value = -1
$ws.Range("B180").Value = 6537869
$ws.Range("C180").Value = "Estonia Meistriliiga"
$ws.Range("D180").Value = "Estonia Meistriliiga"
$ws.Range("E180").Value = 45241.375
$ws.Range("F180").Value = "JK Tallinna Kalev"
$ws.Range("G180").Value = "JK Trans Narva"
$ws.Range("H180").Value = 5
$ws.Range("I180").Value = 0
$ws.Range("J180").Value = "H"
$ws.Range("K180").Value = 1.6
$ws.Range("L180").Value = 4
$ws.Range("M180").Value = 4.5
$ws.Range("N180").Value = 1.65
$ws.Range("O180").Value = 4
$ws.Range("P180").Value = 4.333
$ws.Range("Q180").Value = -0.75
$ws.Range("R180").Value = 1.8
$ws.Range("S180").Value = 2
$ws.Range("T180").Value = 2.75
$ws.Range("U180").Value = 1.9
$ws.Range("V180").Value = 1.9
$ws.Range("W180").Value = 0.6499999999999999
$ws.Range("X180").Value = -1
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = 0.8
$ws.Range("AA180").Value = -1
$ws.Range("AB180").Value = 0.8999999999999999
$ws.Range("AC180").Value = -1
